$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-06-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-29 Saturday", 2) | Out-Null
$d.Content.Find.Execute("484÷6=80, 4", $true, $false, $false, $false, $false, $true, 1, $false, "111÷4=27, 3", 2) | Out-Null
$d.Content.Find.Execute("120÷9=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "192÷4=48, 0", 2) | Out-Null
$d.Content.Find.Execute("930÷6=155, 0", $true, $false, $false, $false, $false, $true, 1, $false, "750÷7=107, 1", 2) | Out-Null
$d.Content.Find.Execute("712÷3=237, 1", $true, $false, $false, $false, $false, $true, 1, $false, "226÷7=32, 2", 2) | Out-Null
$d.Content.Find.Execute("421÷9=46, 7", $true, $false, $false, $false, $false, $true, 1, $false, "405÷7=57, 6", 2) | Out-Null
$d.Content.Find.Execute("113÷4=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "907÷9=100, 7", 2) | Out-Null
$d.Content.Find.Execute("466÷9=51, 7", $true, $false, $false, $false, $false, $true, 1, $false, "845÷9=93, 8", 2) | Out-Null
$d.Content.Find.Execute("855÷7=122, 1", $true, $false, $false, $false, $false, $true, 1, $false, "981÷7=140, 1", 2) | Out-Null
$d.Content.Find.Execute("847÷9=94, 1", $true, $false, $false, $false, $false, $true, 1, $false, "423÷2=211, 1", 2) | Out-Null
$d.Content.Find.Execute("726÷5=145, 1", $true, $false, $false, $false, $false, $true, 1, $false, "476÷5=95, 1", 2) | Out-Null
$d.Content.Find.Execute("855÷9=95, 0", $true, $false, $false, $false, $false, $true, 1, $false, "594÷7=84, 6", 2) | Out-Null
$d.Content.Find.Execute("300÷2=150, 0", $true, $false, $false, $false, $false, $true, 1, $false, "408÷9=45, 3", 2) | Out-Null
$d.Content.Find.Execute("154÷4=38, 2", $true, $false, $false, $false, $false, $true, 1, $false, "199÷8=24, 7", 2) | Out-Null
$d.Content.Find.Execute("750÷6=125, 0", $true, $false, $false, $false, $false, $true, 1, $false, "621÷6=103, 3", 2) | Out-Null
$d.Content.Find.Execute("397÷5=79, 2", $true, $false, $false, $false, $false, $true, 1, $false, "822÷3=274, 0", 2) | Out-Null
$d.Content.Find.Execute("829÷5=165, 4", $true, $false, $false, $false, $false, $true, 1, $false, "344÷2=172, 0", 2) | Out-Null
$d.Content.Find.Execute("474÷8=59, 2", $true, $false, $false, $false, $false, $true, 1, $false, "336÷7=48, 0", 2) | Out-Null
$d.Content.Find.Execute("462÷9=51, 3", $true, $false, $false, $false, $false, $true, 1, $false, "665÷2=332, 1", 2) | Out-Null
$d.Content.Find.Execute("235÷6=39, 1", $true, $false, $false, $false, $false, $true, 1, $false, "286÷5=57, 1", 2) | Out-Null
$d.Content.Find.Execute("800÷5=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "496÷6=82, 4", 2) | Out-Null
$d.Content.Find.Execute("868÷3=289, 1", $true, $false, $false, $false, $false, $true, 1, $false, "240÷4=60, 0", 2) | Out-Null
$d.Content.Find.Execute("487÷9=54, 1", $true, $false, $false, $false, $false, $true, 1, $false, "139÷7=19, 6", 2) | Out-Null
$d.Content.Find.Execute("187÷7=26, 5", $true, $false, $false, $false, $false, $true, 1, $false, "825÷2=412, 1", 2) | Out-Null
$d.Content.Find.Execute("429÷8=53, 5", $true, $false, $false, $false, $false, $true, 1, $false, "959÷2=479, 1", 2) | Out-Null
$d.Content.Find.Execute("471÷9=52, 3", $true, $false, $false, $false, $false, $true, 1, $false, "706÷6=117, 4", 2) | Out-Null
$d.Save()
